# Update "想去人数" (F column) counts on the "展览" sheet and the
# corresponding rows on the "全部类型" sheet (which mirrors the same
# records, offset by one row because of its extra leading entry).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$sheet1Updates = @{
    2  = 142
    3  = 1327
    4  = 1130
    5  = 1019
    6  = 1795
    7  = 561
    8  = 1203
    12 = 295
    13 = 67
    15 = 689
    17 = 103
    18 = 28
    21 = 146
    23 = 34
    24 = 645
    27 = 875
    28 = 315
    31 = 270
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型": same records, shifted down by one row
$sheet4Updates = @{
    3  = 142
    4  = 1327
    5  = 1130
    6  = 1019
    7  = 1795
    8  = 561
    9  = 1203
    14 = 295
    15 = 67
    17 = 689
    19 = 103
    21 = 28
    29 = 146
    31 = 34
    32 = 645
    35 = 875
    36 = 315
    41 = 270
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
